$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (fal6_cropped): update "found" and "correct" columns from [] to [258]
$ws.Range("C3").Value = "[258]"
$ws.Range("D3").Value = "[258]"

# Row 5 (FallingAwayFromCamera): update "found" and "correct" columns from [] to [429]
$ws.Range("C5").Value = "[429]"
$ws.Range("D5").Value = "[429]"
